# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Palta, Vega Modelo de Temuco) above the
# current first data block, pushing the existing rows 1134-1202 down to
# 1137-1205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 1134, shifting everything below down by 3.
$ws.Range("A1134:T1136").EntireRow.Insert()

# --- New row 1134: Especial ---
$ws.Range("A1134").Value = 10
$ws.Range("B1134").Value = "Vega Modelo de Temuco"
$ws.Range("C1134").Value = "La Araucanía"
$ws.Range("D1134").Value = 44753
$ws.Range("E1134").Value = 9
$ws.Range("F1134").Value = "Fruta"
$ws.Range("G1134").Value = 100106
$ws.Range("H1134").Value = "Oleaginosos"
$ws.Range("I1134").Value = 100106002
$ws.Range("J1134").Value = "Palta"
$ws.Range("K1134").Value = "Hass"
$ws.Range("L1134").Value = "Especial"
$ws.Range("M1134").Value = 450
$ws.Range("N1134").Value = 2400
$ws.Range("O1134").Value = 2400
$ws.Range("P1134").Value = 2400
$ws.Range("Q1134").Value = "$/kilo (en caja de 8 kilos )"
$ws.Range("R1134").Value = "Perú"
$ws.Range("S1134").Value = 2400
$ws.Range("T1134").Value = 1

# --- New row 1135: Primera ---
$ws.Range("A1135").Value = 10
$ws.Range("B1135").Value = "Vega Modelo de Temuco"
$ws.Range("C1135").Value = "La Araucanía"
$ws.Range("D1135").Value = 44753
$ws.Range("E1135").Value = 9
$ws.Range("F1135").Value = "Fruta"
$ws.Range("G1135").Value = 100106
$ws.Range("H1135").Value = "Oleaginosos"
$ws.Range("I1135").Value = 100106002
$ws.Range("J1135").Value = "Palta"
$ws.Range("K1135").Value = "Hass"
$ws.Range("L1135").Value = "Primera"
$ws.Range("M1135").Value = 380
$ws.Range("N1135").Value = 1700
$ws.Range("O1135").Value = 1700
$ws.Range("P1135").Value = 1700
$ws.Range("Q1135").Value = "$/kilo (en caja de 8 kilos )"
$ws.Range("R1135").Value = "Perú"
$ws.Range("S1135").Value = 1700
$ws.Range("T1135").Value = 1

# --- New row 1136: Segunda ---
$ws.Range("A1136").Value = 10
$ws.Range("B1136").Value = "Vega Modelo de Temuco"
$ws.Range("C1136").Value = "La Araucanía"
$ws.Range("D1136").Value = 44753
$ws.Range("E1136").Value = 9
$ws.Range("F1136").Value = "Fruta"
$ws.Range("G1136").Value = 100106
$ws.Range("H1136").Value = "Oleaginosos"
$ws.Range("I1136").Value = 100106002
$ws.Range("J1136").Value = "Palta"
$ws.Range("K1136").Value = "Hass"
$ws.Range("L1136").Value = "Segunda"
$ws.Range("M1136").Value = 520
$ws.Range("N1136").Value = 1300
$ws.Range("O1136").Value = 1500
$ws.Range("P1136").Value = 1408
$ws.Range("Q1136").Value = "$/kilo (en caja de 8 kilos )"
$ws.Range("R1136").Value = "Perú"
$ws.Range("S1136").Value = 1408
$ws.Range("T1136").Value = 1
